$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 95283997.94687383
$ws.Range("C2").Value = -181126431.7277304
$ws.Range("I2").Value = 67.83336164715892
$ws.Range("J2").Value = -178.8333050195077
$ws.Range("B3").Value = 119550397.81551
$ws.Range("C3").Value = -209147245.3898863
$ws.Range("I3").Value = 65.69999468594703
$ws.Range("J3").Value = -178.0666719807197
$ws.Range("B4").Value = 104941270.6788485
$ws.Range("C4").Value = -184375531.6367889
$ws.Range("I4").Value = 66.39997687381246
$ws.Range("J4").Value = -174.8833564595209
$ws.Range("B5").Value = 95827846.36066358
$ws.Range("C5").Value = -182845134.7133268
$ws.Range("I5").Value = 67.85001829537772
$ws.Range("J5").Value = -179.2833150379556
$ws.Range("B6").Value = 95741334.71643947
$ws.Range("C6").Value = -182870185.5601847
$ws.Range("I6").Value = 67.86671670282102
$ws.Range("J6").Value = -179.3499499638456
$ws.Range("B7").Value = 110428185.9511601
$ws.Range("C7").Value = -201249734.0073696
$ws.Range("I7").Value = 66.59994315430309
$ws.Range("J7").Value = -179.0500568456969
$ws.Range("B8").Value = 95831097.14160162
$ws.Range("C8").Value = -182596899.4977312
$ws.Range("I8").Value = 67.83337393385574
$ws.Range("J8").Value = -179.1666260661443
$ws.Range("B9").Value = 105424177.0032899
$ws.Range("C9").Value = -186071178.5209873
$ws.Range("I9").Value = 66.43340208612744
$ws.Range("J9").Value = -175.3499312472059
$ws.Range("B10").Value = 93838148.12582983
$ws.Range("C10").Value = -180545283.6552253
$ws.Range("I10").Value = 68.05005468778411
$ws.Range("J10").Value = -179.4999453122159
$ws.Range("B11").Value = 105481876.0653925
$ws.Range("C11").Value = -184961935.7090315
$ws.Range("I11").Value = 66.34997128755195
$ws.Range("J11").Value = -174.8666953791148
$ws.Range("B12").Value = 124228588.4715796
$ws.Range("C12").Value = -206494470.9674475
$ws.Range("I12").Value = 64.86679641086742
$ws.Range("J12").Value = -175.5165369224659
$ws.Range("B13").Value = 114775970.9772273
$ws.Range("C13").Value = -200512542.1981468
$ws.Range("I13").Value = 65.88330810347263
$ws.Range("J13").Value = -176.833358563194
$ws.Range("B14").Value = 92751788.18368557
$ws.Range("C14").Value = -177342903.7849017
$ws.Range("I14").Value = 68.03340471926347
$ws.Range("J14").Value = -178.6999286140698
$ws.Range("B15").Value = 98399684.42395204
$ws.Range("C15").Value = -168147393.1085977
$ws.Range("I15").Value = 66.3832569334397
$ws.Range("J15").Value = -171.1667430665604
$ws.Range("B16").Value = 124690050.8859959
$ws.Range("C16").Value = -205920400.4650612
$ws.Range("I16").Value = 64.76676128743597
$ws.Range("J16").Value = -175.183238712564
$ws.Range("B17").Value = 117110197.7580975
$ws.Range("C17").Value = -202418922.8946863
$ws.Range("I17").Value = 65.64997228439121
$ws.Range("J17").Value = -176.6000277156088
$ws.Range("B18").Value = 121346830.2819801
$ws.Range("C18").Value = -211276373.4495602
$ws.Range("I18").Value = 65.56655332885526
$ws.Range("J18").Value = -178.1334466711447
$ws.Range("B19").Value = 125519215.899418
$ws.Range("C19").Value = -203198296.5584472
$ws.Range("I19").Value = 64.48338593502929
$ws.Range("J19").Value = -174.0499473983041
$ws.Range("B20").Value = 95364596.81383142
$ws.Range("C20").Value = -181344413.6594365
$ws.Range("I20").Value = 67.83346220749966
$ws.Range("J20").Value = -178.883204459167
$ws.Range("B21").Value = 127328664.1379309
$ws.Range("C21").Value = -206191166.5138207
$ws.Range("I21").Value = 64.41671993659487
$ws.Range("J21").Value = -174.4999467300718
$ws.Range("B22").Value = 121266763.0873286
$ws.Range("C22").Value = -198347203.042288
$ws.Range("I22").Value = 64.78343969693921
$ws.Range("J22").Value = -173.7165603030608
$ws.Range("B23").Value = 95381528.85832316
$ws.Range("C23").Value = -182157848.2088327
$ws.Range("I23").Value = 67.88341649871434
$ws.Range("J23").Value = -179.249916834619
$ws.Range("B24").Value = 95605162.06082919
$ws.Range("C24").Value = -182502406.8320286
$ws.Range("I24").Value = 67.86667177744339
$ws.Range("J24").Value = -179.26666155589
$ws.Range("B25").Value = 121302567.4736511
$ws.Range("C25").Value = -203758529.3476131
$ws.Range("I25").Value = 65.11669181224212
$ws.Range("J25").Value = -175.5499748544246
$ws.Range("B26").Value = 112736084.8396494
$ws.Range("C26").Value = -195553994.4118544
$ws.Range("I26").Value = 65.88338446289151
$ws.Range("J26").Value = -175.8166155371085
$ws.Range("B27").Value = 95250469.72453789
$ws.Range("C27").Value = -181546750.717905
$ws.Range("I27").Value = 67.86665547156203
$ws.Range("J27").Value = -179.0500111951046
$ws.Range("B28").Value = 114519157.4111848
$ws.Range("C28").Value = -188018214.6526375
$ws.Range("I28").Value = 65.11670719855645
$ws.Range("J28").Value = -172.3332928014435
$ws.Range("B29").Value = 99552368.42483975
$ws.Range("C29").Value = -173040603.8722652
$ws.Range("I29").Value = 66.53337907415737
$ws.Range("J29").Value = -172.7332875925092
$ws.Range("B30").Value = 122386908.7302104
$ws.Range("C30").Value = -200372427.8239095
$ws.Range("I30").Value = 64.74992377307039
$ws.Range("J30").Value = -174.0500762269296
$ws.Range("B31").Value = 121305639.0987237
$ws.Range("C31").Value = -203762196.873767
$ws.Range("I31").Value = 65.11647441619643
$ws.Range("J31").Value = -175.5501922504702
$ws.Range("B32").Value = 95246610.15562837
$ws.Range("C32").Value = -181792475.6709797
$ws.Range("I32").Value = 67.88331339769263
$ws.Range("J32").Value = -179.1666866023074
$ws.Range("B33").Value = 120390801.7820959
$ws.Range("C33").Value = -211194300.5434513
$ws.Range("I33").Value = 65.69999795399802
$ws.Range("J33").Value = -178.4666687126687
$ws.Range("B34").Value = 95491579.99477631
$ws.Range("C34").Value = -182453648.0181718
$ws.Range("I34").Value = 67.88334847439322
$ws.Range("J34").Value = -179.3166515256068
$ws.Range("B35").Value = 125024109.1397756
$ws.Range("C35").Value = -202349886.662595
$ws.Range("I35").Value = 64.49999526918459
$ws.Range("J35").Value = -173.9166713974821
$ws.Range("B36").Value = 124327141.534732
$ws.Range("C36").Value = -205631350.3833643
$ws.Range("I36").Value = 64.80015036111361
$ws.Range("J36").Value = -175.1998496388864
$ws.Range("B37").Value = 100517524.5837477
$ws.Range("C37").Value = -174700028.5846796
$ws.Range("I37").Value = 66.48322906384377
$ws.Range("J37").Value = -172.9667709361563
$ws.Range("B38").Value = 101868677.8207729
$ws.Range("C38").Value = -175394490.4031445
$ws.Range("I38").Value = 66.29995050572477
$ws.Range("J38").Value = -172.6000494942752
$ws.Range("B39").Value = 92875591.35906237
$ws.Range("C39").Value = -179715338.5746408
$ws.Range("I39").Value = 68.16670171430422
$ws.Range("J39").Value = -179.7499649523625
$ws.Range("B40").Value = 94841119.89169264
$ws.Range("C40").Value = -180446729.0647684
$ws.Range("I40").Value = 67.86673825049411
$ws.Range("J40").Value = -178.7999284161726
$ws.Range("B41").Value = 121229783.9312185
$ws.Range("C41").Value = -203589326.0603937
$ws.Range("I41").Value = 65.11675841779834
$ws.Range("J41").Value = -175.516574915535
$ws.Range("B42").Value = 94348518.45245603
$ws.Range("C42").Value = -180904401.9074258
$ws.Range("I42").Value = 67.98335393134224
$ws.Range("J42").Value = -179.3333127353244
$ws.Range("B43").Value = 93969253.23702408
$ws.Range("C43").Value = -178608145.3456385
$ws.Range("I43").Value = 67.89997591453019
$ws.Range("J43").Value = -178.5000240854698
$ws.Range("B44").Value = 104810639.4110961
$ws.Range("C44").Value = -184303746.9290686
$ws.Range("I44").Value = 66.41679410589634
$ws.Range("J44").Value = -174.9165392274371
$ws.Range("B45").Value = 90477776.63345405
$ws.Range("C45").Value = -173671852.9647405
$ws.Range("I45").Value = 68.20001595428498
$ws.Range("J45").Value = -178.4666507123817
$ws.Range("B46").Value = 95909931.02336974
$ws.Range("C46").Value = -183066432.8012119
$ws.Range("I46").Value = 67.85003020674193
$ws.Range("J46").Value = -179.3333031265913
$ws.Range("B47").Value = 94301473.24302509
$ws.Range("C47").Value = -181548493.462025
$ws.Range("I47").Value = 68.03338138338532
$ws.Range("J47").Value = -179.6666186166147
$ws.Range("B48").Value = 94297103.23677927
$ws.Range("C48").Value = -179236327.3664464
$ws.Range("I48").Value = 67.88328522004672
$ws.Range("J48").Value = -178.5833814466199
$ws.Range("B49").Value = 123511159.385416
$ws.Range("C49").Value = -204562696.9623128
$ws.Range("I49").Value = 64.84995767676973
$ws.Range("J49").Value = -175.1000423232303
$ws.Range("B50").Value = 95726392.43290836
$ws.Range("C50").Value = -181803929.9202596
$ws.Range("I50").Value = 67.80010172042356
$ws.Range("J50").Value = -178.8665649462432
$ws.Range("B51").Value = 105276258.5071207
$ws.Range("C51").Value = -175038955.4961019
$ws.Range("I51").Value = 65.69991906688784
$ws.Range("J51").Value = -170.9167475997789
$ws.Range("B52").Value = 111385813.1348711
$ws.Range("C52").Value = -202325566.3048023
$ws.Range("I52").Value = 66.51671638435909
$ws.Range("J52").Value = -179.0166169489743
$ws.Range("B53").Value = 117762280.2443921
$ws.Range("C53").Value = -206444810.4697082
$ws.Range("I53").Value = 65.80003044384634
$ws.Range("J53").Value = -177.7999695561537
$ws.Range("B54").Value = 90648909.40460366
$ws.Range("C54").Value = -173640852.4905823
$ws.Range("I54").Value = 68.16673432352188
$ws.Range("J54").Value = -178.3332656764781
$ws.Range("B55").Value = 95144979.32683986
$ws.Range("C55").Value = -181007456.3519843
$ws.Range("I55").Value = 67.84998854730327
$ws.Range("J55").Value = -178.8666781193634
$ws.Range("B56").Value = 116765237.1330395
$ws.Range("C56").Value = -201589320.5850689
$ws.Range("I56").Value = 65.65005069036712
$ws.Range("J56").Value = -176.4332826429662
$ws.Range("B57").Value = 90078482.66003077
$ws.Range("C57").Value = -175847742.4345175
$ws.Range("I57").Value = 68.41666202381649
$ws.Range("J57").Value = -179.8166713095169
$ws.Range("B58").Value = 112421186.6435295
$ws.Range("C58").Value = -203309627.5781915
$ws.Range("I58").Value = 66.41674402719615
$ws.Range("J58").Value = -178.9165893061372
$ws.Range("B59").Value = 101353956.1124853
$ws.Range("C59").Value = -174380278.5482855
$ws.Range("I59").Value = 66.31658417360453
$ws.Range("J59").Value = -172.4167491597289
$ws.Range("B60").Value = 108990125.585209
$ws.Range("C60").Value = -198400944.9399034
$ws.Range("I60").Value = 66.65010862310646
$ws.Range("J60").Value = -178.6165580435602
$ws.Range("B61").Value = 122728632.492375
$ws.Range("C61").Value = -218062593.9383368
$ws.Range("I61").Value = 65.76670766974841
$ws.Range("J61").Value = -179.9666256635849
$ws.Range("B62").Value = 123510164.7726939
$ws.Range("C62").Value = -204561546.3511443
$ws.Range("I62").Value = 64.85002801599575
$ws.Range("J62").Value = -175.0999719840042
$ws.Range("B63").Value = 100300642.5028112
$ws.Range("C63").Value = -170669436.9423906
$ws.Range("I63").Value = 66.23326588624691
$ws.Range("J63").Value = -171.3334007804197
$ws.Range("B64").Value = 110373640.5749958
$ws.Range("C64").Value = -200843805.4986471
$ws.Range("I64").Value = 66.58350997521569
$ws.Range("J64").Value = -178.9164900247843
$ws.Range("B65").Value = 102481904.874213
$ws.Range("C65").Value = -192204907.2039707
$ws.Range("I65").Value = 67.31664487291627
$ws.Range("J65").Value = -179.4666884604171
$ws.Range("B66").Value = 105782274.34132
$ws.Range("C66").Value = -174815707.5859639
$ws.Range("I66").Value = 65.60002875733929
$ws.Range("J66").Value = -170.6166379093274
$ws.Range("B67").Value = 111386930.7365307
$ws.Range("C67").Value = -202327117.153652
$ws.Range("I67").Value = 66.51663782919877
$ws.Range("J67").Value = -179.0166955041346
$ws.Range("B68").Value = 99475750.87607253
$ws.Range("C68").Value = -172150960.8091157
$ws.Range("I68").Value = 66.4834134900225
$ws.Range("J68").Value = -172.3832531766441
$ws.Range("B69").Value = 98869816.35752448
$ws.Range("C69").Value = -165438070.2453113
$ws.Range("I69").Value = 66.10010451914354
$ws.Range("J69").Value = -169.7498954808565
$ws.Range("B70").Value = 104219668.5023493
$ws.Range("C70").Value = -177743836.1466576
$ws.Range("I70").Value = 66.06670557210387
$ws.Range("J70").Value = -172.4832944278961
$ws.Range("B71").Value = 117054203.5502286
$ws.Range("C71").Value = -196200399.241827
$ws.Range("I71").Value = 65.26665103512377
$ws.Range("J71").Value = -174.3833489648762
$ws.Range("B72").Value = 95365871.76747249
$ws.Range("C72").Value = -181346368.7595194
$ws.Range("I72").Value = 67.83336632419059
$ws.Range("J72").Value = -178.883300342476
$ws.Range("B73").Value = 121303628.3113402
$ws.Range("C73").Value = -203759795.9989638
$ws.Range("I73").Value = 65.11661673055883
$ws.Range("J73").Value = -175.5500499361079
$ws.Range("B74").Value = 108554960.1844395
$ws.Range("C74").Value = -192823921.6748115
$ws.Range("I74").Value = 66.36665950224786
$ws.Range("J74").Value = -176.6000071644188
$ws.Range("B75").Value = 94542668.27345119
$ws.Range("C75").Value = -181174493.6530331
$ws.Range("I75").Value = 67.96672750150034
$ws.Range("J75").Value = -179.3332724984996
$ws.Range("B76").Value = 126364001.9614239
$ws.Range("C76").Value = -208391483.5858551
$ws.Range("I76").Value = 64.68336175479625
$ws.Range("J76").Value = -175.4833049118704
$ws.Range("B77").Value = 114600057.7617754
$ws.Range("C77").Value = -193023232.9809702
$ws.Range("I77").Value = 65.43335182853889
$ws.Range("J77").Value = -174.1333148381277
$ws.Range("B78").Value = 96182717.2374838
$ws.Range("C78").Value = -167069233.2421187
$ws.Range("I78").Value = 66.70001252623884
$ws.Range("J78").Value = -171.8333208070945
$ws.Range("B79").Value = 90103790.92650808
$ws.Range("C79").Value = -175918167.6152635
$ws.Range("I79").Value = 68.41668237733052
$ws.Range("J79").Value = -179.8333176226694
$ws.Range("B80").Value = 112423939.6619415
$ws.Range("C80").Value = -203313413.7283501
$ws.Range("I80").Value = 66.41655103159917
$ws.Range("J80").Value = -178.9167823017342
$ws.Range("B81").Value = 110880978.5168227
$ws.Range("C81").Value = -201587229.1416554
$ws.Range("I81").Value = 66.54997985485556
$ws.Range("J81").Value = -178.9666868118111
$ws.Range("B82").Value = 111720135.9280964
$ws.Range("C82").Value = -202900426.6929799
$ws.Range("I82").Value = 66.50003988393064
$ws.Range("J82").Value = -179.0832934494026
$ws.Range("B83").Value = 112374935.9972873
$ws.Range("C83").Value = -185557290.790704
$ws.Range("I83").Value = 65.28341023740553
$ws.Range("J83").Value = -172.1999230959278
$ws.Range("B84").Value = 103127393.1932825
$ws.Range("C84").Value = -172088725.1295874
$ws.Range("I84").Value = 65.85007944968972
$ws.Range("J84").Value = -170.6332538836436
$ws.Range("B85").Value = 122710867.0445147
$ws.Range("C85").Value = -203257325.3563437
$ws.Range("I85").Value = 64.8833321995574
$ws.Range("J85").Value = -174.9166678004426
$ws.Range("B86").Value = 111990549.4419555
$ws.Range("C86").Value = -203312206.8927979
$ws.Range("I86").Value = 66.48335103265381
$ws.Range("J86").Value = -179.1166489673462
$ws.Range("B87").Value = 98871311.63442823
$ws.Range("C87").Value = -168604841.0895316
$ws.Range("I87").Value = 66.33342876510508
$ws.Range("J87").Value = -171.1332379015615
$ws.Range("B88").Value = 107111056.3425351
$ws.Range("C88").Value = -195472349.8516857
$ws.Range("I88").Value = 66.76667840048799
$ws.Range("J88").Value = -178.3666549328454
$ws.Range("B89").Value = 90347365.26830955
$ws.Range("C89").Value = -176338122.2630417
$ws.Range("I89").Value = 68.40003475744538
$ws.Range("J89").Value = -179.8666319092213
$ws.Range("B90").Value = 102530379.9016063
$ws.Range("C90").Value = -191802337.6301251
$ws.Range("I90").Value = 67.2832438056665
$ws.Range("J90").Value = -179.2667561943335
$ws.Range("B91").Value = 90103925.83790122
$ws.Range("C91").Value = -175918383.7433024
$ws.Range("I91").Value = 68.41667206489616
$ws.Range("J91").Value = -179.8333279351038
$ws.Range("B92").Value = 99291195.39419353
$ws.Range("C92").Value = -180336735.1781213
$ws.Range("I92").Value = 67.08339802509568
$ws.Range("J92").Value = -176.083268641571
$ws.Range("B93").Value = 89515138.18118136
$ws.Range("C93").Value = -174795692.2836212
$ws.Range("I93").Value = 68.4500570897184
$ws.Range("J93").Value = -179.6999429102816
$ws.Range("B94").Value = 90104599.81563872
$ws.Range("C94").Value = -175919463.4535802
$ws.Range("I94").Value = 68.41662054712215
$ws.Range("J94").Value = -179.8333794528778
$ws.Range("B95").Value = 89566559.09773922
$ws.Range("C95").Value = -174937605.485363
$ws.Range("I95").Value = 68.45001264493334
$ws.Range("J95").Value = -179.7333206884
$ws.Range("B96").Value = 105283032.7366307
$ws.Range("C96").Value = -180788096.2565957
$ws.Range("I96").Value = 66.10008077205472
$ws.Range("J96").Value = -173.2499192279453
$ws.Range("B97").Value = 92664415.74031872
$ws.Range("C97").Value = -179137113.8263955
$ws.Range("I97").Value = 68.16669219156609
$ws.Range("J97").Value = -179.6166411417673
$ws.Range("B98").Value = 105995905.4567777
$ws.Range("C98").Value = -186983498.0103559
$ws.Range("I98").Value = 66.39998833668601
$ws.Range("J98").Value = -175.450011663314
$ws.Range("B99").Value = 111850004.2022698
$ws.Range("C99").Value = -203229053.4905242
$ws.Range("I99").Value = 66.50000423063585
$ws.Range("J99").Value = -179.1499957693642
$ws.Range("B100").Value = 108682215.0231036
$ws.Range("C100").Value = -198149523.846813
$ws.Range("I100").Value = 66.6832680201298
$ws.Range("J100").Value = -178.6667319798702
$ws.Range("B101").Value = 123263194.7934499
$ws.Range("C101").Value = -200773555.3995813
$ws.Range("I101").Value = 64.65007690462897
$ws.Range("J101").Value = -173.9165897620377
$ws.Range("B102").Value = 95284001.11153802
$ws.Range("C102").Value = -181126436.5794339
$ws.Range("I102").Value = 67.83336140894836
$ws.Range("J102").Value = -178.8333052577183
$ws.Range("B103").Value = 103206311.7080425
$ws.Range("C103").Value = -192252878.875446
$ws.Range("I103").Value = 67.19987818182531
$ws.Range("J103").Value = -179.083455151508
$ws.Range("B104").Value = 106044442.6357666
$ws.Range("C104").Value = -192488568.9412396
$ws.Range("I104").Value = 66.74993034355941
$ws.Range("J104").Value = -177.6834029897739
$ws.Range("B105").Value = 97110995.55015253
$ws.Range("C105").Value = -181416775.6881125
$ws.Range("I105").Value = 67.53320201763403
$ws.Range("J105").Value = -177.8334646490326
$ws.Range("B106").Value = 131242611.0262185
$ws.Range("C106").Value = -226691592.6217423
$ws.Range("I106").Value = 65.08346460721658
$ws.Range("J106").Value = -179.83320205945
$ws.Range("B107").Value = 117129515.137243
$ws.Range("C107").Value = -199529806.5510803
$ws.Range("I107").Value = 65.46662074064554
$ws.Range("J107").Value = -175.5500459260211
$ws.Range("B108").Value = 94929093.34019737
$ws.Range("C108").Value = -181961555.8396538
$ws.Range("I108").Value = 67.94991100094209
$ws.Range("J108").Value = -179.4500889990579
$ws.Range("B109").Value = 108291038.0034577
$ws.Range("C109").Value = -179957228.0258081
$ws.Range("I109").Value = 65.54998505796823
$ws.Range("J109").Value = -171.6333482753651
$ws.Range("B110").Value = 88934289.42342946
$ws.Range("C110").Value = -170947381.6615569
$ws.Range("I110").Value = 68.30006331731207
$ws.Range("J110").Value = -178.1999366826879
$ws.Range("B111").Value = 103973242.3097483
$ws.Range("C111").Value = -177149926.3409938
$ws.Range("I111").Value = 66.06655834058715
$ws.Range("J111").Value = -172.3501083260795
$ws.Range("B112").Value = 112619363.6736304
$ws.Range("C112").Value = -204630624.7919917
$ws.Range("I112").Value = 66.46664325015371
$ws.Range("J112").Value = -179.3333567498462
$ws.Range("B113").Value = 90461016.56727831
$ws.Range("C113").Value = -173376641.2661828
$ws.Range("I113").Value = 68.18328501231335
$ws.Range("J113").Value = -178.3333816543533
$ws.Range("B114").Value = 90323452.64385182
$ws.Range("C114").Value = -173250613.5704355
$ws.Range("I114").Value = 68.1999512195265
$ws.Range("J114").Value = -178.3667154471402
$ws.Range("B115").Value = 127185600.9567544
$ws.Range("C115").Value = -222549627.4635976
$ws.Range("I115").Value = 65.39998638549125
$ws.Range("J115").Value = -179.8500136145088
$ws.Range("B116").Value = 104323840.462929
$ws.Range("C116").Value = -179676334.3504549
$ws.Range("I116").Value = 66.18324920806168
$ws.Range("J116").Value = -173.2334174586049
$ws.Range("B117").Value = 104885273.8681174
$ws.Range("C117").Value = -183987818.6619841
$ws.Range("I117").Value = 66.38323079425454
$ws.Range("J117").Value = -174.7501025390788
$ws.Range("B118").Value = 120626242.0701926
$ws.Range("C118").Value = -201107055.5776693
$ws.Range("I118").Value = 65.04992755189618
$ws.Range("J118").Value = -174.8667391147705
$ws.Range("B119").Value = 106254081.1004089
$ws.Range("C119").Value = -191200866.8442713
$ws.Range("I119").Value = 66.63323228473124
$ws.Range("J119").Value = -177.0501010486021
$ws.Range("B120").Value = 93418420.93201588
$ws.Range("C120").Value = -177883323.0553674
$ws.Range("I120").Value = 67.95006425270499
$ws.Range("J120").Value = -178.5166024139617
$ws.Range("B121").Value = 90102621.94079775
$ws.Range("C121").Value = -175916294.8946933
$ws.Range("I121").Value = 68.4167717332719
$ws.Range("J121").Value = -179.833228266728
$ws.Range("B122").Value = 89439994.28489967
$ws.Range("C122").Value = -171832715.3857557
$ws.Range("I122").Value = 68.26657952686253
$ws.Range("J122").Value = -178.2834204731374
$ws.Range("B123").Value = 90056071.40156201
$ws.Range("C123").Value = -172771548.4129456
$ws.Range("I123").Value = 68.21668292665927
$ws.Range("J123").Value = -178.3166504066741
$ws.Range("B124").Value = 112339625.2761361
$ws.Range("C124").Value = -191249031.5810917
$ws.Range("I124").Value = 65.66673111396588
$ws.Range("J124").Value = -174.3499355527008
$ws.Range("B125").Value = 91638675.91084142
$ws.Range("C125").Value = -174580799.3188672
$ws.Range("I125").Value = 68.04990602153245
$ws.Range("J125").Value = -178.1167606451342
$ws.Range("B126").Value = 104042168.0730878
$ws.Range("C126").Value = -182405823.5215728
$ws.Range("I126").Value = 66.41673434655264
$ws.Range("J126").Value = -174.499932320114
$ws.Range("B127").Value = 98433685.40547928
$ws.Range("C127").Value = -164848617.1248325
$ws.Range("I127").Value = 66.13333759743972
$ws.Range("J127").Value = -169.6999957358936
$ws.Range("B128").Value = 90105029.24682693
$ws.Range("C128").Value = -175920151.3993
$ws.Range("I128").Value = 68.41658772219928
$ws.Range("J128").Value = -179.8334122778007
$ws.Range("B129").Value = 105297571.7318893
$ws.Range("C129").Value = -186004414.8481569
$ws.Range("I129").Value = 66.44981829372792
$ws.Range("J129").Value = -175.3835150396054
$ws.Range("B130").Value = 91506452.89152119
$ws.Range("C130").Value = -168656243.1493973
$ws.Range("I130").Value = 67.66673433990181
$ws.Range("J130").Value = -175.3332656600982
$ws.Range("B131").Value = 100244127.4586507
$ws.Range("C131").Value = -181531930.5944287
$ws.Range("I131").Value = 67.00001392915915
$ws.Range("J131").Value = -176.0833194041741
$ws.Range("B132").Value = 90930277.77780448
$ws.Range("C132").Value = -173660515.8562022
$ws.Range("I132").Value = 68.11678675721816
$ws.Range("J132").Value = -178.1498799094485
$ws.Range("B133").Value = 89886613.54807913
$ws.Range("C133").Value = -175570076.8970628
$ws.Range("I133").Value = 68.43329371234678
$ws.Range("J133").Value = -179.8167062876532
$ws.Range("B134").Value = 102207984.6576595
$ws.Range("C134").Value = -187580015.7520474
$ws.Range("I134").Value = 67.06668641223058
$ws.Range("J134").Value = -177.6333135877694
$ws.Range("B135").Value = 105297903.1240131
$ws.Range("C135").Value = -174152576.490499
$ws.Range("I135").Value = 65.63322982578258
$ws.Range("J135").Value = -170.5501035075508
$ws.Range("B136").Value = 94630448.58195198
$ws.Range("C136").Value = -179373449.67
$ws.Range("I136").Value = 67.83336216811897
$ws.Range("J136").Value = -178.4333044985477
$ws.Range("B137").Value = 103648394.1056142
$ws.Range("C137").Value = -182663562.9066028
$ws.Range("I137").Value = 66.49995848417912
$ws.Range("J137").Value = -174.8000415158209
$ws.Range("B138").Value = 98266681.69730829
$ws.Range("C138").Value = -180677399.1216154
$ws.Range("I138").Value = 67.28339650198974
$ws.Range("J138").Value = -176.8166034980103
$ws.Range("B139").Value = 119294187.3056267
$ws.Range("C139").Value = -214136283.2011814
$ws.Range("I139").Value = 66.03337128277319
$ws.Range("J139").Value = -179.9499620505601
$ws.Range("B140").Value = 108799312.6881271
$ws.Range("C140").Value = -196582637.4172918
$ws.Range("I140").Value = 66.56666883258055
$ws.Range("J140").Value = -177.9833311674194
$ws.Range("B141").Value = 90651463.50287198
$ws.Range("C141").Value = -173644848.279627
$ws.Range("I141").Value = 68.1665367111958
$ws.Range("J141").Value = -178.3334632888042
$ws.Range("B142").Value = 109559762.8263735
$ws.Range("C142").Value = -185388356.1909328
$ws.Range("I142").Value = 65.71667090057571
$ws.Range("J142").Value = -173.2333290994243
$ws.Range("B143").Value = 103634673.1074151
$ws.Range("C143").Value = -181891466.9635093
$ws.Range("I143").Value = 66.45001519971581
$ws.Range("J143").Value = -174.4833181336175
$ws.Range("B144").Value = 93414692.50158627
$ws.Range("C144").Value = -163614314.8440734
$ws.Range("I144").Value = 66.94998687113755
$ws.Range("J144").Value = -171.7666797955291
$ws.Range("B145").Value = 119403173.9746964
$ws.Range("C145").Value = -193069939.7414641
$ws.Range("I145").Value = 64.71659064187352
$ws.Range("J145").Value = -172.5000760247931
$ws.Range("B146").Value = 91237591.3195311
$ws.Range("C146").Value = -174741292.29525
$ws.Range("I146").Value = 68.13334226087335
$ws.Range("J146").Value = -178.4666577391267
$ws.Range("B147").Value = 107559905.6503057
$ws.Range("C147").Value = -196354673.0070667
$ws.Range("I147").Value = 66.7500683016478
$ws.Range("J147").Value = -178.4999316983522
$ws.Range("B148").Value = 125912277.5563406
$ws.Range("C148").Value = -202194389.5591151
$ws.Range("I148").Value = 64.36657390169819
$ws.Range("J148").Value = -173.6167594316352
$ws.Range("B149").Value = 107665495.2339031
$ws.Range("C149").Value = -194506578.9795287
$ws.Range("I149").Value = 66.61661649160766
$ws.Range("J149").Value = -177.700050175059
$ws.Range("B150").Value = 96673180.40597731
$ws.Range("C150").Value = -182291982.2595517
$ws.Range("I150").Value = 67.66685949910348
$ws.Range("J150").Value = -178.4998071675632
$ws.Range("B151").Value = 105035929.1928367
$ws.Range("C151").Value = -185858192.3171963
$ws.Range("I151").Value = 66.48325698771632
$ws.Range("J151").Value = -175.450076345617
$ws.Range("B152").Value = 129710198.4596708
$ws.Range("C152").Value = -224508383.7818565
$ws.Range("I152").Value = 65.16663848931606
$ws.Range("J152").Value = -179.6333615106839
$ws.Range("B153").Value = 110114046.0511149
$ws.Range("C153").Value = -198837048.2266916
$ws.Range("I153").Value = 66.49995514339913
$ws.Range("J153").Value = -178.2500448566009
$ws.Range("B154").Value = 95345677.07905038
$ws.Range("C154").Value = -180780445.0931121
$ws.Range("I154").Value = 67.7998927347974
$ws.Range("J154").Value = -178.6334405985359
$ws.Range("B155").Value = 119325984.7731324
$ws.Range("C155").Value = -200735864.3579961
$ws.Range("I155").Value = 65.21671717554545
$ws.Range("J155").Value = -175.1832828244545
$ws.Range("B156").Value = 93428580.87039356
$ws.Range("C156").Value = -177405798.8254772
$ws.Range("I156").Value = 67.91659312823128
$ws.Range("J156").Value = -178.2834068717687
$ws.Range("B157").Value = 114373753.1311606
$ws.Range("C157").Value = -190705215.0722503
$ws.Range("I157").Value = 65.31665241504226
$ws.Range("J157").Value = -173.3666809182911
$ws.Range("B158").Value = 104586924.6323144
$ws.Range("C158").Value = -174359046.694427
$ws.Range("I158").Value = 65.76664134229371
$ws.Range("J158").Value = -170.9333586577063
$ws.Range("B159").Value = 91749763.98484507
$ws.Range("C159").Value = -175631218.9138206
$ws.Range("I159").Value = 68.09999197765791
$ws.Range("J159").Value = -178.5500080223421
$ws.Range("B160").Value = 125836086.652237
$ws.Range("C160").Value = -202024684.8876894
$ws.Range("I160").Value = 64.36662366929453
$ws.Range("J160").Value = -173.5833763307054
$ws.Range("B161").Value = 116302363.5055453
$ws.Range("C161").Value = -190096358.4436287
$ws.Range("I161").Value = 64.98335676951956
$ws.Range("J161").Value = -172.4666432304805
$ws.Range("B162").Value = 107310720.9184249
$ws.Range("C162").Value = -195718577.6450581
$ws.Range("I162").Value = 66.7500911044866
$ws.Range("J162").Value = -178.3665755621801
$ws.Range("B163").Value = 114563223.8634599
$ws.Range("C163").Value = -192168517.7273502
$ws.Range("I163").Value = 65.3833839032541
$ws.Range("J163").Value = -173.8332827634125
$ws.Range("B164").Value = 95010238.99019727
$ws.Range("C164").Value = -180392876.9759969
$ws.Range("I164").Value = 67.83343937331193
$ws.Range("J164").Value = -178.6665606266881
$ws.Range("B165").Value = 116834089.474993
$ws.Range("C165").Value = -190551582.2354816
$ws.Range("I165").Value = 64.93333182453715
$ws.Range("J165").Value = -172.4500015087962
$ws.Range("B166").Value = 123867543.5620008
$ws.Range("C166").Value = -203490825.668832
$ws.Range("I166").Value = 64.73336601050909
$ws.Range("J166").Value = -174.6333006561576
$ws.Range("B167").Value = 94930724.12062098
$ws.Range("C167").Value = -180177071.7016331
$ws.Range("I167").Value = 67.83327419627541
$ws.Range("J167").Value = -178.6167258037246
$ws.Range("B168").Value = 109023935.8231804
$ws.Range("C168").Value = -190379709.2240946
$ws.Range("I168").Value = 66.13328115770848
$ws.Range("J168").Value = -175.4167188422916
$ws.Range("B169").Value = 122611482.311808
$ws.Range("C169").Value = -195357198.7251206
$ws.Range("I169").Value = 64.40003391967613
$ws.Range("J169").Value = -172.3166327469906
$ws.Range("B170").Value = 112012002.4043519
$ws.Range("C170").Value = -187947364.5359129
$ws.Range("I170").Value = 65.50014057935047
$ws.Range("J170").Value = -173.2331927539828
$ws.Range("B171").Value = 115464219.5656805
$ws.Range("C171").Value = -189436241.734163
$ws.Range("I171").Value = 65.0666289865151
$ws.Range("J171").Value = -172.5167043468183
$ws.Range("B172").Value = 123165732.0728984
$ws.Range("C172").Value = -204036426.8111035
$ws.Range("I172").Value = 64.86662275926916
$ws.Range("J172").Value = -175.0333772407308
$ws.Range("B173").Value = 126036634.0646268
$ws.Range("C173").Value = -201936231.7772607
$ws.Range("I173").Value = 64.33330987357901
$ws.Range("J173").Value = -173.5000234597543
$ws.Range("B174").Value = 103962297.5280395
$ws.Range("C174").Value = -184188421.1273362
$ws.Range("I174").Value = 66.55001094406209
$ws.Range("J174").Value = -175.2833223892712
$ws.Range("B175").Value = 123167497.8179095
$ws.Range("C175").Value = -196595138.8558727
$ws.Range("I175").Value = 64.39998428865744
$ws.Range("J175").Value = -172.5666823780093
$ws.Range("B176").Value = 125007944.8906292
$ws.Range("C176").Value = -206926604.33199
$ws.Range("I176").Value = 64.78352735916813
$ws.Range("J176").Value = -175.4164726408319
$ws.Range("B177").Value = 104730402.1973267
$ws.Range("C177").Value = -183606540.2147338
$ws.Range("I177").Value = 66.38326247672906
$ws.Range("J177").Value = -174.666737523271
$ws.Range("B178").Value = 108429037.3583531
$ws.Range("C178").Value = -197506227.4677483
$ws.Range("I178").Value = 66.68339334081354
$ws.Range("J178").Value = -178.5332733258531
$ws.Range("B179").Value = 107512582.6285103
$ws.Range("C179").Value = -195967449.1250856
$ws.Range("I179").Value = 66.73335430959045
$ws.Range("J179").Value = -178.3666456904096
$ws.Range("B180").Value = 119548694.4906461
$ws.Range("C180").Value = -208863939.1255731
$ws.Range("I180").Value = 65.68320561056693
$ws.Range("J180").Value = -177.9667943894331
